$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(3.0468275264548565, 1.239625383183412, -0.014177555829265314, 2.618972214012707)"
$ws.Range("C2").Value = "JSU(-1.2128937733613094, 1.3339925296175839, 2.2894003183464853, 5.908483215763431)"
$ws.Range("D2").Value = "NIG(0.7201114721969599, 0.4595527151399929, 1.7099238987612904, 2.797237825461088)"
$ws.Range("E2").Value = "JSU(-1.543523806209776, 1.3778128622180574, 1.544314991092791, 5.372569775371369)"
